$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "保險" (insurance) - 5th worksheet
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Value = "company"
$ws5.Range("C1").Value = "name"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "property_category"
$ws5.Range("F1").Value = "category"
$ws5.Range("G1").Value = "date"
$ws5.Range("H1").Value = "legislator_name"
$ws5.Range("I1").Value = "legislator_id"
$ws5.Range("J1").Value = "source_file"
$ws5.Range("K1").Value = "index"

$ws5.Range("B2").Value = "富邦人壽"
$ws5.Range("C2").Value = "安泰人壽55年增值分紅養老壽險"
$ws5.Range("D2").Value = "詹文馨"
$ws5.Range("E2").Value = "insurance"
$ws5.Range("F2").Value = "normal"
$ws5.Range("G2").Value = "2012-04-09"
$ws5.Range("H2").Value = "吳秉叡"
$ws5.Range("I2").Value = 1324
$ws5.Range("J2").Value = "tmp8f8d1"
$ws5.Range("K2").Value = 107

# ---------------------------------------------------------------------------
# Sheet "債務" (debt) - 6th worksheet
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Value = "species"
$ws6.Range("C1").Value = "debtor"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"
$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

$ws6.Range("B2").Value = "抵押借款"
$ws6.Range("C2").Value = "吳秉叡"
$ws6.Range("D2").Value = "大台北商業銀行城內分行臺北市中正區重慶南路"
$ws6.Range("E2").Value = 9580000
$ws6.Range("F2").Value = "101年02月16日"
$ws6.Range("G2").Value = "投資"
$ws6.Range("H2").Value = "debt"
$ws6.Range("I2").Value = "normal"
$ws6.Range("J2").Value = "2012-04-09"
$ws6.Range("K2").Value = "吳秉叡"
$ws6.Range("L2").Value = 1324
$ws6.Range("M2").Value = "tmp8f8d1"
$ws6.Range("N2").Value = 117

$ws6.Range("B3").Value = "抵押借款"
$ws6.Range("C3").Value = "吳秉叡"
$ws6.Range("D3").Value = "大台北商業銀行城内分行臺北市中正區重慶南路"
$ws6.Range("E3").Value = 500000
$ws6.Range("F3").Value = "101年03月14曰"
$ws6.Range("G3").Value = "投資"
$ws6.Range("H3").Value = "debt"
$ws6.Range("I3").Value = "normal"
$ws6.Range("J3").Value = "2012-04-09"
$ws6.Range("K3").Value = "吳秉叡"
$ws6.Range("L3").Value = 1324
$ws6.Range("M3").Value = "tmp8f8d1"
$ws6.Range("N3").Value = 118

# ---------------------------------------------------------------------------
# Sheet "事業投資" (investment) - 7th worksheet
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

$ws7.Range("B1").Value = "owner"
$ws7.Range("C1").Value = "company"
$ws7.Range("D1").Value = "address"
$ws7.Range("E1").Value = "total"
$ws7.Range("F1").Value = "register_date"
$ws7.Range("G1").Value = "register_reason"
$ws7.Range("H1").Value = "property_category"
$ws7.Range("I1").Value = "category"
$ws7.Range("J1").Value = "date"
$ws7.Range("K1").Value = "legislator_name"
$ws7.Range("L1").Value = "legislator_id"
$ws7.Range("M1").Value = "source_file"
$ws7.Range("N1").Value = "index"

$ws7.Range("B2").Value = "吳秉數"
$ws7.Range("C2").Value = "雄裕建設股份有限公司"
$ws7.Range("D2").Value = "臺北市忠孝東路五段508號14樓之3"
$ws7.Range("E2").Value = 4372000
$ws7.Range("F2").Value = "100年03月11曰"
$ws7.Range("G2").Value = "個人投資"
$ws7.Range("H2").Value = "investment"
$ws7.Range("I2").Value = "normal"
$ws7.Range("J2").Value = "2012-04-09"
$ws7.Range("K2").Value = "吳秉叡"
$ws7.Range("L2").Value = 1324
$ws7.Range("M2").Value = "tmp8f8d1"
$ws7.Range("N2").Value = 123

$ws7.Range("B3").Value = "吳秉叡"
$ws7.Range("C3").Value = "聖裕投資有限公司"
$ws7.Range("D3").Value = "新北市龍安路106巷1號4樓"
$ws7.Range("E3").Value = 4801500
$ws7.Range("F3").Value = "100年08月29日"
$ws7.Range("G3").Value = "個人投資"
$ws7.Range("H3").Value = "investment"
$ws7.Range("I3").Value = "normal"
$ws7.Range("J3").Value = "2012-04-09"
$ws7.Range("K3").Value = "吳秉叡"
$ws7.Range("L3").Value = 1324
$ws7.Range("M3").Value = "tmp8f8d1"
$ws7.Range("N3").Value = 124

$ws7.Range("B4").Value = "吳秉叡"
$ws7.Range("C4").Value = "益翔建設股份有限公司"
$ws7.Range("D4").Value = "新北市成泰路一段2號3樓"
$ws7.Range("E4").Value = 10000000
$ws7.Range("F4").Value = "99年02月23日"
$ws7.Range("G4").Value = "個人投資"
$ws7.Range("H4").Value = "investment"
$ws7.Range("I4").Value = "normal"
$ws7.Range("J4").Value = "2012-04-09"
$ws7.Range("K4").Value = "吳秉叡"
$ws7.Range("L4").Value = 1324
$ws7.Range("M4").Value = "tmp8f8d1"
$ws7.Range("N4").Value = 125

$ws7.Range("B5").Value = "吳秉叡"
$ws7.Range("C5").Value = "福益實業股份有限公司"
$ws7.Range("D5").Value = "新北市成泰路一段2號3樓"
$ws7.Range("E5").Value = 1650
$ws7.Range("F5").Value = "78年07月01H"
$ws7.Range("G5").Value = "個人投資"
$ws7.Range("H5").Value = "investment"
$ws7.Range("I5").Value = "normal"
$ws7.Range("J5").Value = "2012-04-09"
$ws7.Range("K5").Value = "吳秉叡"
$ws7.Range("L5").Value = 1324
$ws7.Range("M5").Value = "tmp8f8d1"
$ws7.Range("N5").Value = 126
